# Add new aural observation rows 489-511 (2020 is complete)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A489").Value = "great crested flycatcher"
$ws.Range("B489").Value = "birds"
$ws.Range("C489").Value = 44006
$ws.Range("D489").Value = 0.77083333333333337
$ws.Range("F489").Value = "Maple Grove"

$ws.Range("A490").Value = "american toad"
$ws.Range("B490").Value = "herps"
$ws.Range("C490").Value = 44006
$ws.Range("D490").Value = 0.78819444444444453
$ws.Range("F490").Value = "Maple Grove"

$ws.Range("A491").Value = "eastern wood peewee"
$ws.Range("B491").Value = "birds"
$ws.Range("C491").Value = 44006
$ws.Range("D491").Value = 0.78819444444444453
$ws.Range("F491").Value = "Maple Grove"

$ws.Range("A492").Value = "wood thrush"
$ws.Range("B492").Value = "birds"
$ws.Range("C492").Value = 44006
$ws.Range("D492").Value = 0.78819444444444453
$ws.Range("F492").Value = "Maple Grove"

$ws.Range("A493").Value = "eastern wood peewee"
$ws.Range("B493").Value = "birds"
$ws.Range("C493").Value = 43979
$ws.Range("D493").Value = 0.83333333333333337
$ws.Range("F493").Value = "Maple Grove"
$ws.Range("K493").Value = "time unspecified, but after the news"

$ws.Range("A494").Value = "great crested flycatcher"
$ws.Range("B494").Value = "birds"
$ws.Range("C494").Value = 43980
$ws.Range("D494").Value = 0.30208333333333331
$ws.Range("F494").Value = "Maple Grove"

$ws.Range("A495").Value = "red-eyed vireo"
$ws.Range("B495").Value = "birds"
$ws.Range("C495").Value = 43980
$ws.Range("D495").Value = 0.30208333333333331
$ws.Range("F495").Value = "Maple Grove"

$ws.Range("A496").Value = "cicada"
$ws.Range("B496").Value = "insects"
$ws.Range("C496").Value = 43984
$ws.Range("D496").Value = 0.625
$ws.Range("F496").Value = "Otis"
$ws.Range("K496").Value = "time unspecified, but hot afternoon"

$ws.Range("A497").Value = "cicada"
$ws.Range("B497").Value = "insects"
$ws.Range("C497").Value = 43994
$ws.Range("D497").Value = 0.29166666666666669
$ws.Range("F497").Value = "Otis"
$ws.Range("K497").Value = "bracketing… `"By a few days later, the trees purred with cicadas, and as I write this on the 12th of June, they still do. Each morning, Brooklyn and I walk through town beneath a cloud of murmurs and humming, as the cicadas gear up for the noisy day.`""

$ws.Range("A498").Value = "cicada"
$ws.Range("B498").Value = "insects"
$ws.Range("C498").Value = 44024
$ws.Range("D498").Value = 0.32291666666666669
$ws.Range("F498").Value = "Maple Grove"

$ws.Range("A499").Value = "great horned owl"
$ws.Range("B499").Value = "birds"
$ws.Range("C499").Value = 44128
$ws.Range("D499").Value = 0.75
$ws.Range("F499").Value = "Maple Grove"

$ws.Range("A500").Value = "golden-crowned kinglet"
$ws.Range("B500").Value = "birds"
$ws.Range("C500").Value = 44128
$ws.Range("D500").Value = 0.75
$ws.Range("F500").Value = "Maple Grove"

$ws.Range("A501").Value = "robin"
$ws.Range("B501").Value = "birds"
$ws.Range("C501").Value = 44049
$ws.Range("D501").Value = 0.41666666666666669
$ws.Range("F501").Value = "Morton Arboretum"
$ws.Range("K501").Value = "time unspecified, but morning moss work with Wayne"

$ws.Range("A502").Value = "eastern wood peewee"
$ws.Range("B502").Value = "birds"
$ws.Range("C502").Value = 44049
$ws.Range("D502").Value = 0.41666666666666669
$ws.Range("F502").Value = "Morton Arboretum"
$ws.Range("K502").Value = "time unspecified, but morning moss work with Wayne"

$ws.Range("A503").Value = "goldfinch"
$ws.Range("B503").Value = "birds"
$ws.Range("C503").Value = 44050
$ws.Range("D503").Value = 0.3125
$ws.Range("F503").Value = "Maple Grove"
$ws.Range("K503").Value = "time from field notebook"

$ws.Range("A504").Value = "red-eyed vireo"
$ws.Range("B504").Value = "birds"
$ws.Range("C504").Value = 44050
$ws.Range("D504").Value = 0.3125
$ws.Range("F504").Value = "Maple Grove"
$ws.Range("K504").Value = "time from field notebook"

$ws.Range("A505").Value = "blue jay"
$ws.Range("B505").Value = "birds"
$ws.Range("C505").Value = 44050
$ws.Range("D505").Value = 0.3125
$ws.Range("F505").Value = "Maple Grove"
$ws.Range("K505").Value = "time from field notebook"

$ws.Range("A506").Value = "robin"
$ws.Range("B506").Value = "birds"
$ws.Range("C506").Value = 44056
$ws.Range("D506").Value = 0.27083333333333331
$ws.Range("F506").Value = "Maple Grove"

$ws.Range("A507").Value = "great horned owl"
$ws.Range("B507").Value = "birds"
$ws.Range("C507").Value = 44056
$ws.Range("D507").Value = 0.27083333333333331
$ws.Range("F507").Value = "Maple Grove"

$ws.Range("A508").Value = "chorus frog"
$ws.Range("B508").Value = "birds"
$ws.Range("C508").Value = 44144
$ws.Range("D508").Value = 0.625
$ws.Range("F508").Value = "Morton Arboretum"
$ws.Range("K508").Value = "times estimated"

$ws.Range("A509").Value = "bluebird"
$ws.Range("B509").Value = "birds"
$ws.Range("C509").Value = 44146
$ws.Range("D509").Value = 0.29166666666666669
$ws.Range("F509").Value = "Otis"
$ws.Range("K509").Value = "times estimated"

$ws.Range("A510").Value = "junco"
$ws.Range("B510").Value = "birds"
$ws.Range("C510").Value = 44146
$ws.Range("D510").Value = 0.29166666666666669
$ws.Range("F510").Value = "Otis"
$ws.Range("K510").Value = "times estimated"

$ws.Range("A511").Value = "sandhill crane"
$ws.Range("B511").Value = "birds"
$ws.Range("C511").Value = 44156
$ws.Range("D511").Value = 0.5
$ws.Range("F511").Value = "Otis"
$ws.Range("K511").Value = "purely a placeholder: the cranes streamed by all day"

[void]$ws.Range("A504").Select()

Write-Host "done"